# Swap the match-data (columns B, and E through AB) between pairs of rows.
# Columns A (id), C (Div) and D (Date) stay attached to their own row.
# NOTE: use Value2 (not Value) - Value has issues reflecting numeric
# cells correctly in this COM-interop runtime.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs whose "data" columns (B, E:AB) get swapped with each other.
$pairs = @(
    @(26, 27),
    @(38, 39),
    @(43, 44),
    @(47, 48),
    @(54, 55)
)

# Columns to swap: B (2), and E..AB (5..28)
$cols = @(2) + (5..28)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)

        $val1 = $cell1.Value2
        $val2 = $cell2.Value2

        $cell1.Value2 = $val2
        $cell2.Value2 = $val1
    }
}
